# Applies the "plotting" sheet addition + regulators/source plotting flags +
# notes "Location" column + misc selection changes described in the commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "plotting" worksheet right after "master".
# ---------------------------------------------------------------------------
$master = $wb.Worksheets.Item("master")
$plotting = $wb.Worksheets.Add($null, $master)
$plotting.Name = "plotting"

# Header row (row 2). Write in this order so new shared strings land at the
# same indices the authored workbook uses (106..110).
$plotting.Range("D2").Value = "Regulators"
$plotting.Range("C2").Value = "Source"
$plotting.Range("D7").Value = "na"
$plotting.Range("E2").Value = "BUS ID"
$plotting.Range("F2").Value = "SRC Reg"

# Row 3 - 13 bus
$plotting.Range("A3").Value = "IEEE Circuits"
$plotting.Range("B3").Value = "13 bus"
$plotting.Range("C3").Value = "X"
$plotting.Range("D3").Value = "X"

# Row 4 - 34 bus
$plotting.Range("B4").Value = "34 bus"
$plotting.Range("C4").Value = "X"
$plotting.Range("D4").Value = "X"

# Row 5 - 123 bus
$plotting.Range("B5").Value = "123 bus"
$plotting.Range("C5").Value = "X"
$plotting.Range("D5").Value = "X"
$plotting.Range("F5").Value = "Y"

# Row 6 - EU LV
$plotting.Range("B6").Value = "EU LV"

# Row 7 - EPRI test
$plotting.Range("B7").Value = "EPRI test"
$plotting.Range("C7").Value = "X"
$plotting.Range("D7").Value = "na"
$plotting.Range("E7").Value = 1

# Row 8 - Ckt5
$plotting.Range("A8").Value = "EPRI test"
$plotting.Range("B8").Value = "Ckt5"
$plotting.Range("C8").Value = "X"
$plotting.Range("D8").Value = "na"
$plotting.Range("E8").Value = 796

# Row 9 - Ckt7
$plotting.Range("B9").Value = "Ckt7"
$plotting.Range("C9").Value = "X"
$plotting.Range("D9").Value = "na"
$plotting.Range("E9").Value = 318412

# Row 10 - Ckt24
$plotting.Range("B10").Value = "Ckt24"

# Row 11 - 8500 node
$plotting.Range("A11").Value = "EPRI HC"
$plotting.Range("B11").Value = "8500 node"
$plotting.Range("C11").Value = "X"
$plotting.Range("D11").Value = 0
$plotting.Range("E11").Value = 5964927408

# Row 12 - V2 (9000 LFs)
$plotting.Range("B12").Value = "V2 (9000 LFs)"
$plotting.Range("C12").Value = "X"
$plotting.Range("D12").Value = "X"
$plotting.Range("F12").Value = "Y"

# Row 13 - 6000 LFs
$plotting.Range("B13").Value = "6000 LFs"
$plotting.Range("C13").Value = "X"
$plotting.Range("D13").Value = "X"
$plotting.Range("F13").Value = "Y"

# Formatting touch-ups to mirror the other sheets' group-header borders.
$master.Range("A3").Copy()
$plotting.Range("A3").PasteSpecial(-4122)
$plotting.Range("A8").PasteSpecial(-4122)
$plotting.Range("A11").PasteSpecial(-4122)

$master.Range("B3").Copy()
$plotting.Range("B8").PasteSpecial(-4122)
$plotting.Range("B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$plotting.Columns.Item(4).ColumnWidth = 10.43

$plotting.Range("P9").Select()

# ---------------------------------------------------------------------------
# 2) "notes" sheet - new "Location" column (J) for the circuits used in the
#    regulator/source plots.
# ---------------------------------------------------------------------------
$notes = $wb.Worksheets.Item("notes")

$notes.Range("J1").Value = "Location"
$notes.Range("J14").Value = "SouthEast"
$notes.Range("J4").Value = "NorthEast"
$notes.Range("J5").Value = "NorthWest"
$notes.Range("J7").Value = "NorthEast"
$notes.Range("J9").Value = "NorthEast"
$notes.Range("J11").Value = "NorthWest"
$notes.Range("J12").Value = "NorthWest"
$notes.Range("J15").Value = "NorthWest"
$notes.Range("J16").Value = "NorthWest"

$notes.Columns.Item(10).ColumnWidth = 10.71

# ---------------------------------------------------------------------------
# 3) Selection / active-sheet bookkeeping to match the saved view state.
# ---------------------------------------------------------------------------
$master.Activate()
$master.Range("A3:B16").Select()

$varStuff = $wb.Worksheets.Item("varStuff")
$varStuff.Activate()
$varStuff.Range("L21").Select()

$notes.Activate()
$notes.Range("H7").Select()
